{"js": "// Replace the date line and each multiplication-problem cell's text\n// with its updated value, per the commit's regenerated numbers.\nconst replacements = [\n  [\"2025-02-10 Monday\", \"2025-02-11 Tuesday\"],\n  [\"749\u00d75=\", \"544\u00d79=\"],\n  [\"316\u00d74=\", \"159\u00d72=\"],\n  [\"340\u00d78=\", \"664\u00d78=\"],\n  [\"813\u00d79=\", \"717\u00d77=\"],\n  [\"357\u00d74=\", \"211\u00d73=\"],\n  [\"345\u00d73=\", \"383\u00d72=\"],\n  [\"863\u00d79=\", \"185\u00d75=\"],\n  [\"264\u00d72=\", \"642\u00d75=\"],\n  [\"858\u00d78=\", \"605\u00d76=\"],\n  [\"411\u00d79=\", \"227\u00d72=\"],\n  [\"330\u00d78=\", \"620\u00d79=\"],\n  [\"176\u00d74=\", \"217\u00d73=\"],\n  [\"874\u00d72=\", \"313\u00d75=\"],\n  [\"956\u00d76=\", \"327\u00d77=\"],\n  [\"170\u00d79=\", \"633\u00d73=\"],\n  [\"102\u00d78=\", \"617\u00d77=\"],\n  [\"870\u00d73=\", \"508\u00d74=\"],\n  [\"440\u00d79=\", \"341\u00d75=\"],\n  [\"417\u00d73=\", \"315\u00d73=\"],\n  [\"113\u00d79=\", \"793\u00d72=\"],\n  [\"353\u00d76=\", \"414\u00d75=\"],\n  [\"819\u00d77=\", \"324\u00d75=\"],\n  [\"229\u00d72=\", \"434\u00d77=\"],\n  [\"629\u00d78=\", \"710\u00d75=\"],\n  [\"186\u00d72=\", \"217\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each multiplication-problem cell's text\n# with its updated value, per the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-10 Monday\", \"2025-02-11 Tuesday\"),\n    @(\"749\u00d75=\", \"544\u00d79=\"),\n    @(\"316\u00d74=\", \"159\u00d72=\"),\n    @(\"340\u00d78=\", \"664\u00d78=\"),\n    @(\"813\u00d79=\", \"717\u00d77=\"),\n    @(\"357\u00d74=\", \"211\u00d73=\"),\n    @(\"345\u00d73=\", \"383\u00d72=\"),\n    @(\"863\u00d79=\", \"185\u00d75=\"),\n    @(\"264\u00d72=\", \"642\u00d75=\"),\n    @(\"858\u00d78=\", \"605\u00d76=\"),\n    @(\"411\u00d79=\", \"227\u00d72=\"),\n    @(\"330\u00d78=\", \"620\u00d79=\"),\n    @(\"176\u00d74=\", \"217\u00d73=\"),\n    @(\"874\u00d72=\", \"313\u00d75=\"),\n    @(\"956\u00d76=\", \"327\u00d77=\"),\n    @(\"170\u00d79=\", \"633\u00d73=\"),\n    @(\"102\u00d78=\", \"617\u00d77=\"),\n    @(\"870\u00d73=\", \"508\u00d74=\"),\n    @(\"440\u00d79=\", \"341\u00d75=\"),\n    @(\"417\u00d73=\", \"315\u00d73=\"),\n    @(\"113\u00d79=\", \"793\u00d72=\"),\n    @(\"353\u00d76=\", \"414\u00d75=\"),\n    @(\"819\u00d77=\", \"324\u00d75=\"),\n    @(\"229\u00d72=\", \"434\u00d77=\"),\n    @(\"629\u00d78=\", \"710\u00d75=\"),\n    @(\"186\u00d72=\", \"217\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
